$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Alexander"
$ws.Range("C2").Value = "Volkov"
$ws.Range("D2").Value = 0

$ws.Range("B3").Value = "Jairzinho"
$ws.Range("C3").Value = "Rozenstruik"
$ws.Range("D3").Value = 0

$ws.Range("B4").Value = "Dan"
$ws.Range("C4").Value = "Ige"
$ws.Range("D4").Value = 0

$ws.Range("B5").Value = "Movsar"
$ws.Range("C5").Value = "Evloev"
$ws.Range("D5").Value = 1

$ws.Range("B6").Value = "Michael"
$ws.Range("C6").Value = "Trizano"
$ws.Range("D6").Value = 1

$ws.Range("B7").Value = "Lucas"
$ws.Range("C7").Value = "Almeida"
$ws.Range("D7").Value = 1

$ws.Range("B8").Value = "Karine"
$ws.Range("C8").Value = "Silva"
$ws.Range("D8").Value = 0

$ws.Range("B9").Value = "Poliana"
$ws.Range("C9").Value = "Botelho"
$ws.Range("D9").Value = 0

$ws.Range("B10").Value = "Ode"
$ws.Range("C10").Value = "Osbourne"
$ws.Range("D10").Value = 0

$ws.Range("B11").Value = "Zarrukh"
$ws.Range("C11").Value = "Adashev"
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = "Alonzo"
$ws.Range("C12").Value = "Menifield"
$ws.Range("D12").Value = 0

$ws.Range("B13").Value = "Askar"
$ws.Range("C13").Value = "Mozharov"
$ws.Range("D13").Value = 1

$ws.Range("B14").Value = "Felice"
$ws.Range("C14").Value = "Herrig"
$ws.Range("D14").Value = 0

$ws.Range("B15").Value = "Karolina"
$ws.Range("C15").Value = "Kowalkiewicz"
$ws.Range("D15").Value = 0

$ws.Range("B16").Value = "Joe"
$ws.Range("C16").Value = "Solecki"
$ws.Range("D16").Value = 0

$ws.Range("B17").Value = "Alex"
$ws.Range("C17").Value = "Da Silva"
$ws.Range("D17").Value = 0

$ws.Range("B18").Value = "Damon"
$ws.Range("C18").Value = "Jackson"
$ws.Range("D18").Value = 0

$ws.Range("B19").Value = "Dan"
$ws.Range("C19").Value = "Argueta"
$ws.Range("D19").Value = 0

$ws.Range("B20").Value = "Niklas"
$ws.Range("C20").Value = "Stolze"
$ws.Range("D20").Value = 0

$ws.Range("B21").Value = "Benoit"
$ws.Range("C21").Value = "Saint Denis"
$ws.Range("D21").Value = 0

$ws.Range("B22").Value = "Johnny"
$ws.Range("C22").Value = "Munoz"
$ws.Range("D22").Value = 0

$ws.Range("B23").Value = "Tony"
$ws.Range("C23").Value = "Gravely"
$ws.Range("D23").Value = 0

$ws.Range("B24").Value = "Jeff"
$ws.Range("C24").Value = "Molina"
$ws.Range("D24").Value = 0

$ws.Range("B25").Value = "Zhalgas"
$ws.Range("C25").Value = "Zhumagulov"
$ws.Range("D25").Value = 0

$ws.Range("B26").Value = "Rinat"
$ws.Range("C26").Value = "Fakhretdinov"
$ws.Range("D26").Value = 0

$ws.Range("B27").Value = "Andreas"
$ws.Range("C27").Value = "Michailidis"
$ws.Range("D27").Value = 0

$ws.Range("B28").Value = "Erin"
$ws.Range("C28").Value = "Blanchfield"
$ws.Range("D28").Value = 0

$ws.Range("B29").Value = "JJ"
$ws.Range("C29").Value = "Aldrich"
$ws.Range("D29").Value = 0
